$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C2:C18) from 2023-09-06 to 2023-09-14
# (serial date 45175 -> 45183), preserving existing date formatting.
for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 3).Value = 45183
}
